# Fix typo in Ch 4 exercise plot dates (column F, rows 2-7): add 15 days.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 44974
$ws.Range("F3").Value = 44973
$ws.Range("F4").Value = 44972
$ws.Range("F5").Value = 44971
$ws.Range("F6").Value = 44970
$ws.Range("F7").Value = 44969
